{"js": "\n// 1. Bookmark \"smartrewards\": renumber id from 1 -> 0.\n// This document only has a single bookmark, and bookmark ids are\n// re-assigned by Word/engine on save; deleting + re-adding it at the\n// same (collapsed) range causes the id to be recomputed starting at 0.\n{\n  const doc = context.document;\n  const bmRange = doc.getBookmarkRange(\"smartrewards\");\n  doc.deleteBookmark(\"smartrewards\");\n  bmRange.insertBookmark(\"smartrewards\");\n  await context.sync();\n}\n\n// 2. Paragraph pPr change + run split for the \"SmartRewards are a price\n// stabilization...\" paragraph. We replace the whole paragraph (pPr +\n// runs) with fully specified OOXML so every property lands exactly.\n{\n  const body = context.document.body;\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"text\");\n  await context.sync();\n\n  let target = null;\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text.indexOf(\"are a price stabilization mechanism\") !== -1) {\n      target = paragraphs.items[i];\n      break;\n    }\n  }\n  if (!target) {\n    throw new Error(\"Target paragraph not found\");\n  }\n\n  const newParaOoxml = \"<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\\\"Times New Roman\\\" w:eastAsia=\\\"Times New Roman\\\" w:hAnsi=\\\"Times New Roman\\\" w:cs=\\\"Times New Roman\\\"/><w:sz w:val=\\\"24\\\"/><w:szCs w:val=\\\"24\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\\\"inherit\\\" w:hAnsi=\\\"inherit\\\" w:cs=\\\"Times New Roman\\\" w:eastAsia=\\\"Times New Roman\\\"/><w:b/><w:bCs/><w:color w:val=\\\"2B2B2B\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/><w:bdr w:val=\\\"none\\\" w:sz=\\\"0\\\" w:space=\\\"0\\\" w:color=\\\"auto\\\" w:frame=\\\"1\\\"/></w:rPr><w:t>SmartRewards</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Times New Roman\\\" w:eastAsia=\\\"Times New Roman\\\"/><w:color w:val=\\\"252525\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\">&#x00a0;are a price stabilization mechanism and a way to encourage long term holding. Long term holders are key to the project\\u2019s success since the SmartHive treasury needs SmartCash to appreciate in value in order to fund meaningful 3rd party proposals and help grow SmartCash into a successful global crypto-currency. Beginning at block 574,100, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Open Sans\\\" w:hAnsi=\\\"Open Sans\\\" w:cs=\\\"Open Sans\\\" w:eastAsia=\\\"Times New Roman\\\"/><w:color w:val=\\\"252525\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/><w:shd w:val=\\\"clear\\\" w:color=\\\"auto\\\" w:fill=\\\"FFFFFF\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\">the current SmartRewards model changed to a Decentralized Distribution</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Times New Roman\\\" w:eastAsia=\\\"Times New Roman\\\"/><w:color w:val=\\\"252525\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\">. SmartRewards will then be distributed after every&#x00a0;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"inherit\\\" w:hAnsi=\\\"inherit\\\" w:cs=\\\"Times New Roman\\\" w:eastAsia=\\\"Times New Roman\\\"/><w:b/><w:bCs/><w:color w:val=\\\"2B2B2B\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/><w:bdr w:val=\\\"none\\\" w:sz=\\\"0\\\" w:space=\\\"0\\\" w:color=\\\"auto\\\" w:frame=\\\"1\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\">47,500 Blocks</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Times New Roman\\\" w:eastAsia=\\\"Times New Roman\\\"/><w:color w:val=\\\"252525\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\">&#x00a0;to all eligible addresses. The payouts will begin 200 Blocks after the cycle ends and 1000 addresses will be paid every other block. SmartNodes will keep track on which addresses are eligible and which are not. You will earn SmartRewards on&#x00a0;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"inherit\\\" w:hAnsi=\\\"inherit\\\" w:cs=\\\"Times New Roman\\\" w:eastAsia=\\\"Times New Roman\\\"/><w:b/><w:bCs/><w:color w:val=\\\"2B2B2B\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/><w:bdr w:val=\\\"none\\\" w:sz=\\\"0\\\" w:space=\\\"0\\\" w:color=\\\"auto\\\" w:frame=\\\"1\\\"/></w:rPr><w:t>any</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Times New Roman\\\" w:eastAsia=\\\"Times New Roman\\\"/><w:color w:val=\\\"252525\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\">&#x00a0;address for which you hold the keys (web or desktop, SmartNode included!) which holds &gt;=1000 SMART for one month and does not make&#x00a0;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"inherit\\\" w:hAnsi=\\\"inherit\\\" w:cs=\\\"Times New Roman\\\" w:eastAsia=\\\"Times New Roman\\\"/><w:b/><w:bCs/><w:color w:val=\\\"2B2B2B\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/><w:bdr w:val=\\\"none\\\" w:sz=\\\"0\\\" w:space=\\\"0\\\" w:color=\\\"auto\\\" w:frame=\\\"1\\\"/></w:rPr><w:t>any</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Times New Roman\\\" w:eastAsia=\\\"Times New Roman\\\"/><w:color w:val=\\\"252525\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\">outgoing transactions during that time. Please note, most exchanges do&#x00a0;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"inherit\\\" w:hAnsi=\\\"inherit\\\" w:cs=\\\"Times New Roman\\\" w:eastAsia=\\\"Times New Roman\\\"/><w:i/><w:iCs/><w:color w:val=\\\"252525\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/><w:bdr w:val=\\\"none\\\" w:sz=\\\"0\\\" w:space=\\\"0\\\" w:color=\\\"auto\\\" w:frame=\\\"1\\\"/></w:rPr><w:t>not</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Times New Roman\\\" w:eastAsia=\\\"Times New Roman\\\"/><w:color w:val=\\\"252525\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\">&#x00a0;pay SmartRewards to their users, holding &gt;= 1000 SMART on an exchange does not guarantee a reward</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"inherit\\\" w:hAnsi=\\\"inherit\\\" w:cs=\\\"Times New Roman\\\" w:eastAsia=\\\"Times New Roman\\\"/><w:b/><w:bCs/><w:color w:val=\\\"2B2B2B\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/><w:bdr w:val=\\\"none\\\" w:sz=\\\"0\\\" w:space=\\\"0\\\" w:color=\\\"auto\\\" w:frame=\\\"1\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\">.&#x00a0;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Times New Roman\\\" w:eastAsia=\\\"Times New Roman\\\"/><w:color w:val=\\\"252525\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\">The SmartRewards will come out of the 15% block reward allocation.</w:t></w:r></w:p>\";\n\n  const pkg =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + newParaOoxml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n\n  target.insertOoxml(pkg, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word COM interop script: reproduce the OOXML diff for\n# \"WHAT IS SMARTCASH.docx\" (smartrewards bookmark id + SmartRewards\n# paragraph formatting/text update).\n\n$d = $word.ActiveDocument\n\n# 1. Bookmark \"smartrewards\": renumber id from 1 -> 0.\n# This document only has a single bookmark; bookmark ids are re-assigned\n# on save, so deleting + re-adding the bookmark at the same (collapsed)\n# location causes its id to be recomputed starting at 0.\n$bm = $d.Bookmarks.Item(\"smartrewards\")\n$bmRange = $bm.Range\n$bm.Delete()\n$d.Bookmarks.Add(\"smartrewards\", $bmRange)\n\n# 2. Paragraph pPr change + run split for the \"SmartRewards are a price\n# stabilization...\" paragraph. Replace the whole paragraph (pPr + runs)\n# with fully specified OOXML so every property lands exactly as in the\n# target revision.\n$paras = $d.Paragraphs\n$target = $null\nfor ($i = 1; $i -le $paras.Count; $i++) {\n  $p = $paras.Item($i)\n  if ($p.Range.Text -like \"*are a price stabilization mechanism*\") {\n    $target = $p\n    break\n  }\n}\nif ($target -eq $null) {\n  throw \"Target paragraph not found\"\n}\n\n$newParaOoxml = '<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:eastAsia=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"inherit\" w:hAnsi=\"inherit\" w:cs=\"Times New Roman\" w:eastAsia=\"Times New Roman\"/><w:b/><w:bCs/><w:color w:val=\"2B2B2B\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/><w:bdr w:val=\"none\" w:sz=\"0\" w:space=\"0\" w:color=\"auto\" w:frame=\"1\"/></w:rPr><w:t>SmartRewards</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Times New Roman\" w:eastAsia=\"Times New Roman\"/><w:color w:val=\"252525\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/></w:rPr><w:t xml:space=\"preserve\">&#x00a0;are a price stabilization mechanism and a way to encourage long term holding. Long term holders are key to the project\u2019s success since the SmartHive treasury needs SmartCash to appreciate in value in order to fund meaningful 3rd party proposals and help grow SmartCash into a successful global crypto-currency. Beginning at block 574,100, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Open Sans\" w:hAnsi=\"Open Sans\" w:cs=\"Open Sans\" w:eastAsia=\"Times New Roman\"/><w:color w:val=\"252525\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\">the current SmartRewards model changed to a Decentralized Distribution</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Times New Roman\" w:eastAsia=\"Times New Roman\"/><w:color w:val=\"252525\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/></w:rPr><w:t xml:space=\"preserve\">. SmartRewards will then be distributed after every&#x00a0;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"inherit\" w:hAnsi=\"inherit\" w:cs=\"Times New Roman\" w:eastAsia=\"Times New Roman\"/><w:b/><w:bCs/><w:color w:val=\"2B2B2B\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/><w:bdr w:val=\"none\" w:sz=\"0\" w:space=\"0\" w:color=\"auto\" w:frame=\"1\"/></w:rPr><w:t xml:space=\"preserve\">47,500 Blocks</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Times New Roman\" w:eastAsia=\"Times New Roman\"/><w:color w:val=\"252525\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/></w:rPr><w:t xml:space=\"preserve\">&#x00a0;to all eligible addresses. The payouts will begin 200 Blocks after the cycle ends and 1000 addresses will be paid every other block. SmartNodes will keep track on which addresses are eligible and which are not. You will earn SmartRewards on&#x00a0;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"inherit\" w:hAnsi=\"inherit\" w:cs=\"Times New Roman\" w:eastAsia=\"Times New Roman\"/><w:b/><w:bCs/><w:color w:val=\"2B2B2B\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/><w:bdr w:val=\"none\" w:sz=\"0\" w:space=\"0\" w:color=\"auto\" w:frame=\"1\"/></w:rPr><w:t>any</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Times New Roman\" w:eastAsia=\"Times New Roman\"/><w:color w:val=\"252525\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/></w:rPr><w:t xml:space=\"preserve\">&#x00a0;address for which you hold the keys (web or desktop, SmartNode included!) which holds &gt;=1000 SMART for one month and does not make&#x00a0;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"inherit\" w:hAnsi=\"inherit\" w:cs=\"Times New Roman\" w:eastAsia=\"Times New Roman\"/><w:b/><w:bCs/><w:color w:val=\"2B2B2B\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/><w:bdr w:val=\"none\" w:sz=\"0\" w:space=\"0\" w:color=\"auto\" w:frame=\"1\"/></w:rPr><w:t>any</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Times New Roman\" w:eastAsia=\"Times New Roman\"/><w:color w:val=\"252525\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/></w:rPr><w:t xml:space=\"preserve\">outgoing transactions during that time. Please note, most exchanges do&#x00a0;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"inherit\" w:hAnsi=\"inherit\" w:cs=\"Times New Roman\" w:eastAsia=\"Times New Roman\"/><w:i/><w:iCs/><w:color w:val=\"252525\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/><w:bdr w:val=\"none\" w:sz=\"0\" w:space=\"0\" w:color=\"auto\" w:frame=\"1\"/></w:rPr><w:t>not</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Times New Roman\" w:eastAsia=\"Times New Roman\"/><w:color w:val=\"252525\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/></w:rPr><w:t xml:space=\"preserve\">&#x00a0;pay SmartRewards to their users, holding &gt;= 1000 SMART on an exchange does not guarantee a reward</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"inherit\" w:hAnsi=\"inherit\" w:cs=\"Times New Roman\" w:eastAsia=\"Times New Roman\"/><w:b/><w:bCs/><w:color w:val=\"2B2B2B\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/><w:bdr w:val=\"none\" w:sz=\"0\" w:space=\"0\" w:color=\"auto\" w:frame=\"1\"/></w:rPr><w:t xml:space=\"preserve\">.&#x00a0;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Times New Roman\" w:eastAsia=\"Times New Roman\"/><w:color w:val=\"252525\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/></w:rPr><w:t xml:space=\"preserve\">The SmartRewards will come out of the 15% block reward allocation.</w:t></w:r></w:p>'\n\n$pkg = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + $newParaOoxml + '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>'\n\n$target.Range.InsertXML($pkg)\n"}
